# Update yearly income-statement data: drop oldest period (1396/12),
# shift remaining periods left, and add the newly reported period (1401/12)
# with its recomputed trailing-12-month figures. Also refresh the
# publish-date header row and fix the previously blank "-" cell to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(8, 4).Value = '12 ماهه منتهی به 1397/12'
$ws.Cells.Item(8, 5).Value = '12 ماهه منتهی به 1398/12'
$ws.Cells.Item(8, 6).Value = '12 ماهه منتهی به 1399/12'
$ws.Cells.Item(8, 7).Value = '12 ماهه منتهی به 1400/12'
$ws.Cells.Item(8, 8).Value = '12 ماهه منتهی به 1401/12'
$ws.Cells.Item(9, 4).Value = '1399-01-27 (10)'
$ws.Cells.Item(9, 5).Value = '1400-02-05 (8)'
$ws.Cells.Item(9, 6).Value = '1401-02-10 (9)'
$ws.Cells.Item(9, 7).Value = '1402-02-12 (10)'
$ws.Cells.Item(9, 8).Value = '1402-02-12 (2)'
$ws.Cells.Item(11, 4).Value = 1753984
$ws.Cells.Item(11, 5).Value = 2730411
$ws.Cells.Item(11, 6).Value = 4546249
$ws.Cells.Item(11, 7).Value = 8166123
$ws.Cells.Item(11, 8).Value = 12146108
$ws.Cells.Item(12, 4).Value = -1282039
$ws.Cells.Item(12, 5).Value = -1877173
$ws.Cells.Item(12, 6).Value = -2990949
$ws.Cells.Item(12, 7).Value = -4982668
$ws.Cells.Item(12, 8).Value = -7538475
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 853238
$ws.Cells.Item(13, 6).Value = 1555300
$ws.Cells.Item(13, 7).Value = 3183455
$ws.Cells.Item(13, 8).Value = 4607633
$ws.Cells.Item(14, 4).Value = -168472
$ws.Cells.Item(14, 5).Value = -302047
$ws.Cells.Item(14, 6).Value = -347085
$ws.Cells.Item(14, 7).Value = -439067
$ws.Cells.Item(14, 8).Value = -744022
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(16, 4).Value = -75074
$ws.Cells.Item(16, 5).Value = -71550
$ws.Cells.Item(16, 6).Value = 82158
$ws.Cells.Item(16, 7).Value = -344346
$ws.Cells.Item(16, 8).Value = 113993
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 479641
$ws.Cells.Item(17, 6).Value = 1290373
$ws.Cells.Item(17, 7).Value = 2400042
$ws.Cells.Item(17, 8).Value = 3977604
$ws.Cells.Item(18, 4).Value = -82765
$ws.Cells.Item(18, 5).Value = -66869
$ws.Cells.Item(18, 6).Value = -14871
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = -7373
$ws.Cells.Item(19, 4).Value = -11102
$ws.Cells.Item(19, 5).Value = 37535
$ws.Cells.Item(19, 6).Value = 353829
$ws.Cells.Item(19, 7).Value = 237498
$ws.Cells.Item(19, 8).Value = 492256
$ws.Cells.Item(20, 4).Value = 0
$ws.Cells.Item(20, 5).Value = 450307
$ws.Cells.Item(20, 6).Value = 1629331
$ws.Cells.Item(20, 7).Value = 2637540
$ws.Cells.Item(20, 8).Value = 4462487
$ws.Cells.Item(21, 4).Value = -17440
$ws.Cells.Item(21, 5).Value = -78000
$ws.Cells.Item(21, 6).Value = -163874
$ws.Cells.Item(21, 7).Value = -341243
$ws.Cells.Item(21, 8).Value = -409181
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(22, 5).Value = 372307
$ws.Cells.Item(22, 6).Value = 1465457
$ws.Cells.Item(22, 7).Value = 2296297
$ws.Cells.Item(22, 8).Value = 4053306
$ws.Cells.Item(23, 4).Value = 0
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(24, 4).Value = 0
$ws.Cells.Item(24, 5).Value = 372307
$ws.Cells.Item(24, 6).Value = 1465457
$ws.Cells.Item(24, 7).Value = 2296297
$ws.Cells.Item(24, 8).Value = 4053306
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 338
$ws.Cells.Item(25, 6).Value = 1332
$ws.Cells.Item(25, 7).Value = 2088
$ws.Cells.Item(25, 8).Value = 3685
$ws.Cells.Item(26, 4).Value = 550000
$ws.Cells.Item(26, 5).Value = 1100000
$ws.Cells.Item(26, 6).Value = 1100000
$ws.Cells.Item(26, 7).Value = 1100000
$ws.Cells.Item(26, 8).Value = 1100000
$ws.Cells.Item(27, 4).Value = 0
$ws.Cells.Item(27, 5).Value = 338
$ws.Cells.Item(27, 6).Value = 1332
$ws.Cells.Item(27, 7).Value = 2088
$ws.Cells.Item(27, 8).Value = 3685
